$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 22:22"

# --- Update numeric stats for several existing countries ---

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 912010
$ws.Cells.Item(4, 3).Value = 25568
$ws.Cells.Item(4, 5).Value = 768291
$ws.Cells.Item(4, 7).Value = 1217
$ws.Cells.Item(4, 8).Value = 51453

# Row 14: Brasil
$ws.Cells.Item(14, 2).Value = 52995
$ws.Cells.Item(14, 3).Value = 3503
$ws.Cells.Item(14, 5).Value = 22752
$ws.Cells.Item(14, 7).Value = 357
$ws.Cells.Item(14, 8).Value = 3670

# Row 16: Canada
$ws.Cells.Item(16, 2).Value = 43552
$ws.Cells.Item(16, 3).Value = 1442
$ws.Cells.Item(16, 4).Value = 15469
$ws.Cells.Item(16, 5).Value = 25789

# Row 122: Isla de Man
$ws.Cells.Item(122, 2).Value = 308
$ws.Cells.Item(122, 3).Value = 1
$ws.Cells.Item(122, 4).Value = 230
$ws.Cells.Item(122, 5).Value = 60
$ws.Cells.Item(122, 7).Value = 2
$ws.Cells.Item(122, 8).Value = 18

# Row 148: Monaco
$ws.Cells.Item(148, 4).Value = 41
$ws.Cells.Item(148, 5).Value = 49

# --- Re-order "Guyana" ahead of "San Martin (Parte Holandesa)" / Haiti / Bahamas,
#     and refresh its case numbers to the latest figures ---

# Insert a fresh row just above "San Martin (Parte Holandesa)" (row 156) and push
# San Martin / Haiti / Bahamas / Islas Caimanes down by one row.
$ws.Rows.Item(156).Insert()

$ws.Cells.Item(156, 1).Value = "Guyana"
$ws.Cells.Item(156, 2).Value = 73
$ws.Cells.Item(156, 3).Value = 3
$ws.Cells.Item(156, 4).Value = 12
$ws.Cells.Item(156, 5).Value = 54
$ws.Cells.Item(156, 6).Value = 5
$ws.Cells.Item(156, 7).Value = 0
$ws.Cells.Item(156, 8).Value = 7

# Remove the old "Guyana" row, which has now shifted down to row 160.
$ws.Rows.Item(160).Delete()
